# "Added periodic & upfront related scenarios"
#
# The ProductLoanInput sheet's "repaymentstrategy" row (B17) previously held
# the placeholder value "Mifos style". Populate it with the real strategy
# description and give it the left/top-aligned, wrap-friendly look used for
# the other descriptive answer cells (same green input fill, just aligned to
# the top-left instead of centered), then make this sheet the one that's
# active/selected when the workbook is reopened, with the cursor resting on
# the cell that was just filled in.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

$target = $ws1.Range("B17")
$target.Value = "Penalties, Fees, Interest, Principal order"

# Left/top alignment on the same green "input" fill already used elsewhere.
$target.HorizontalAlignment = -4131   # xlLeft
$target.VerticalAlignment = -4160     # xlTop
$target.Interior.Color = 5296274      # RGB(146,208,80) == existing fillId 2

# ProductLoanInput becomes the active/selected sheet (was ProductLoanOutput).
$target.Select()
$ws1.Activate()
